# Rename the "GRA-vehbatsubsidy" worksheet to "GRA-batconsubsidy".
# (Battery production subsidy cash-flow subscript rename.)
$wb = $excel.ActiveWorkbook
$sheet = $wb.Worksheets.Item("GRA-vehbatsubsidy")
$sheet.Name = "GRA-batconsubsidy"
